# Auto-generated edit script: updates cryptos list cell values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'67.789.32"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +7.75%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'3.551.38"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +10.21%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.01%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'192.21"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +10.29%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'560.97"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +9.00%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'3.548.19"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +10.20%  "

$ws.Range("E8").Value = "  +3.49%  "

$ws.Range("E9").Value = "  -0.07%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.643"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +7.60%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'56.52"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +7.79%  "

$ws.Range("E12").Value = "  +16.49%  "

$ws.Range("E13").Value = "  +9.43%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'9.52"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +7.36%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'4.112.01"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +10.17%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.545.15"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +10.28%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'67.787.75"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +7.93%  "

$ws.Range("E18").Value = "  +5.63%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'18.46"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +7.80%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'11.96"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +9.27%  "

$ws.Range("E21").Value = "  +5.08%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'407.66"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +11.66%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'85.57"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +7.11%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'4.25"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +9.64%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'11.46"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +4.24%  "

$ws.Range("E27").Value = "  +15.02%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'6.16"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("E29").Value = "  +7.24%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'8.89"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +9.12%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'30.69"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +8.90%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'692.79"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +6.39%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'6.86"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +9.35%  "

$ws.Range("E34").Value = "  +6.96%  "

$ws.Range("E35").Value = "  +8.74%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'60.92"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +5.76%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.0₃0834"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +19.93%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'39.19"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +7.43%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  +7.25%  "

$ws.Range("E41").Value = "  +15.37%  "

$ws.Range("E42").Value = "  +18.39%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'3.067.72"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +7.64%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'3.00"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +16.39%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.70"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +7.55%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'3.33"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +13.93%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.0423"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +9.06%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'9.11"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +20.26%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'2.74"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +2.29%  "

$ws.Range("E51").Value = "  +7.22%  "
